# Applies the diff described in the commit:
#   "changed MP time limit and corrected error in fixed recourse data"
#
# Two kinds of changes to the workbook:
#  1) "Sheet1" summary sheet (rows 2..11 / instances 1..10): the fixed
#     recourse data was corrected - new objective/solve-time values and
#     new num-cuts / num-variables / num-cons / num-quad-cons counts
#     (reflecting the larger MP time limit => more cuts generated).
#  2) Each per-instance detail sheet ("1".."10"): row 2 (iteration 1) gets
#     a new MP solve time / worst violation, and row 3 (iteration 2) gets
#     a new MP objective / MP gap / MP solve time (consequence of the
#     longer MP time limit run).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Summary sheet ("Sheet1"), columns B,C,F,G,H,I for rows 2..11 ---
$summary = @(
  @{ Row=2;  B=-1446.9652710430485; C=815.022106596;  F=100; G=110400; H=120500; I=10000 },
  @{ Row=3;  B=-1446.960829055565;  C=729.902236155;   F=100; G=110400; H=120500; I=10000 },
  @{ Row=4;  B=-1447.5007260736406; C=529.191510378;   F=100; G=110400; H=120500; I=10000 },
  @{ Row=5;  B=-1443.4704065646756; C=892.866737132;   F=100; G=110400; H=120500; I=10000 },
  @{ Row=6;  B=-1452.0331302759485; C=981.336032032;   F=100; G=110400; H=120500; I=10000 },
  @{ Row=7;  B=-1434.9010419212857; C=3607.487382638;  F=100; G=110400; H=120500; I=10000 },
  @{ Row=8;  B=-1448.9218603408262; C=1000.859303818;  F=100; G=110400; H=120500; I=10000 },
  @{ Row=9;  B=-1434.8107321419998; C=832.827194559;   F=100; G=110400; H=120500; I=10000 },
  @{ Row=10; B=-1443.5569295182886; C=631.220203737;   F=100; G=110400; H=120500; I=10000 },
  @{ Row=11; B=-1453.0961120243596; C=1049.538362124;  F=100; G=110400; H=120500; I=10000 }
)

foreach ($r in $summary) {
  $ws1.Range("B" + $r.Row).Value = $r.B
  $ws1.Range("C" + $r.Row).Value = $r.C
  $ws1.Range("F" + $r.Row).Value = $r.F
  $ws1.Range("G" + $r.Row).Value = $r.G
  $ws1.Range("H" + $r.Row).Value = $r.H
  $ws1.Range("I" + $r.Row).Value = $r.I
}

# --- Per-instance detail sheets ("1".."10") ---
# row2: D (MP solve time), E (Worst violation)
# row3: B (MP objective), C (MP gap), D (MP solve time)
$details = @(
  @{ Sheet="9";  D2=0.10140063171020508;  E2=166.29438; B3=-1443.5569295182886; C3=0.07613152052995122; D3=623.0101736463762 },
  @{ Sheet="10"; D2=0.0684218231298828;   E2=161.68958; B3=-1453.0961120243596; C3=0.09219493660558657; D3=1041.5924777124642 },
  @{ Sheet="1";  D2=0.9683720088411866;   E2=163.46956; B3=-1446.9652710430485; C3=0.08898265372079207; D3=799.7161838207189 },
  @{ Sheet="2";  D2=0.1123262791986084;   E2=164.5962;  B3=-1446.960829055565;  C3=0.01993502586008528; D3=720.8830310219234 },
  @{ Sheet="3";  D2=0.06738471238122559;  E2=163.91369; B3=-1447.5007260736406; C3=0.09715401022847754; D3=523.0575813672768 },
  @{ Sheet="4";  D2=0.08500793670581054;  E2=165.41896; B3=-1443.4704065646756; C3=0.08205533729814819; D3=884.9945875280558 },
  @{ Sheet="5";  D2=0.1028845121800537;   E2=166.67449; B3=-1452.0331302759485; C3=0.09787459685418172; D3=974.975940314295 },
  @{ Sheet="6";  D2=0.1104670791965332;   E2=164.15669; B3=-1434.9010419212857; C3=0.14142294036018835; D3=3600.6305427022016 },
  @{ Sheet="7";  D2=0.06842382542602539;  E2=164.13388; B3=-1448.9218603408262; C3=0.07599605688155607; D3=994.9392517533083 },
  @{ Sheet="8";  D2=0.10054509087524414;  E2=167.05476; B3=-1434.8107321419998; C3=0.09931419573023853; D3=825.1823572381986 }
)

foreach ($d in $details) {
  $ws = $wb.Worksheets.Item($d.Sheet)
  $ws.Range("D2").Value = $d.D2
  $ws.Range("E2").Value = $d.E2
  $ws.Range("B3").Value = $d.B3
  $ws.Range("C3").Value = $d.C3
  $ws.Range("D3").Value = $d.D3
}
